$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet "My Series" -> "Data"
$ws.Name = "Data"

# 2. Update cell A11 text
$ws.Range("A11").Value = "Function Information"

# 3. Update cell B21 value (Kurtosis)
$ws.Range("B21").Value = 0.2499825759175085

# 4. Update the number format used by B27:B36 (custom numFmt 166: "0.000" -> "###0.000")
$ws.Range("B27:B36").NumberFormat = "###0.000"

# 5. Update the comment attached to A1 with the refreshed embedded metadata blob
$comment = $ws.Range("A1").Comment
[void]$comment.Text("Jh0AAB+LCAAAAAAAAAOlWVtvI0kV/istP4GE3d12ZjaJKr3yLVkLO45sh0z2BZW7K3GRdrfpqk7it0UCLVoWIYRm0XIVT4uQGEawKy0zXP7LapIZnvgLnLr0zXZ20mEUTbrO+U7VqVPnVhX07vXcNy5JxGgY7FXsmlUxSOCGHg3O9yoxP6vajyvvOqh77RL/CEd4TjiADZAK2O41o3uVGeeLXdO8urqqXTVqYXRu1i3LNp8M+mN3Rua4SgPGceCSSirlvV2q4qC2Nx8Qjj3MsZLcq/TGvVqbULcDtAEO8DmJaq2Y0YAw1g045ZQwIRkRzEm7M/iO2phTrz2u2chco2fIVkx9T+EKSEXXOFiWTOicOHXL3q5a29WGNbHru1s7u/Wdmt2w308EUyDqY8bHJLqkriSMOZ4vpLi13bDs+tZOfQeZG0EwV2YABw19b0QuKSNem/g+K2URUx9g0+Ww63LGtJCZk9UTPVyFgwgvZhPKfVJOjdGgZcwDrUs2iYP2w4i4YL8HqXRIroaRNutk0QfuZEYjvuzgZem5jhmJhgthpHKiDuqEAW/6JOLHCzhr4oErAMPhUUyQeQczE+pQ5sI3DWLiOWfYZ3mhAhOdhNEFW2CXHEIcm2KOq8APsQcOxynj1M0WXWOgoyhcwIyweCv0vX2YVYM3MNKZewGYWCzbCsOLTLtNTCRPVZ4vnOkc8wS+RkfjWXg1DPzlOJ4yN6JT4nVaCXojD4mA1NLtmPFwDlpkJKRoOcoS/kEArpJRh7h0jv0jH4zInAbMUiCgZszDM8rboR/PA5botEJFJ7CjCblODJgO0RDONhA2D4NekMAVaCOrKDAKr9Il1xnSBjlyk7nJaa8zVsEdoCWnt86RByI2uU+9KA/5o8hRi04xnhHCN3qE4iCRCfdFwXFay8N4PoXomkKIXcpVGTIzPgI3BVcHvRwLikhV/kwsa1f+gB4pG3UD725cwkSwXG4txwbeCgnBnvyWj4MLoJ5QPjtsJnvZwEHKAnfi13kIAnfh46Ukp1bK01AvcP3YIyof9IIz6aFCN3Wod7LRGqkPIe4gHCwnywWkZUZ3OXzsVaBQ7zIeQStQcdwwDni0FIkDmRr6NhkWTwO5APbvLXMWke/H0IEs9+PAbYfe/VfzlHWOA8rvr2EYRyob3l9EWk8kxph1iEgxMuffW94tsycWlYLPAzIPA+re39pgZKG994CNsCSq7i1BVHzdG+9DWVdVT8T6vcUiaB+hzpVapslY6FLprDo8vJy8eUfIdMgZjn1o3ThU2PM0966SUZNdrGLyJHQc+UkGdERjzKAzdr15zYXeQXR/NTecC4IJDenJGJl5vGiAXNINzvs4OI+hxUjzyio9zb+iPE4iHDCxnbSjWEnFm0EoyVOq03FU8hrG0hFU8gqBi8wVHJqQ+SKMsD8Aw9B97Xa6XYJOZID5TI+gtvnETYxsZqKpVFGzRPG3wWSRUtsQAa/T5ApRgsReVA+eYTIaErscQFj6bezTaaSyqlZiIw8OLOsNk/wrNleyT0zOAO5hUH2/TZaiOc8Gmi5d1k4YyoFFInXGo63t+iOrUYfGRoyR3PGIYN/oQjBzYvSCS8L4HMR2jRFh1IMviv1d4z0yJRSKoDSRLkOlpfNyaD/J81KVJuhbpBQB0G+cUygj68CUkwk4pwRH/jIHVFvthy7gbn/875vfPH/14rPbj56++eKH//3Hr17982c3z34EH7d//dvNx79U21RgNMFTn0iFJq3tbauxBX6WkpAwrik7Yy92uaSdnsqGOB0jfa+Tg3a31z7ot2Q+SYmJuCopprgyLsM4G47VJuRC8kjNxBMUxJkk+UmPC9xciXLEFe6SFNF5/l2CyhavX372+uWf75TWBst6LXtn51HVrr+1FYPLsL2GS1uxfqEGCPBW1XpUrddz4BUMGqkCkNqp5zkN29qx6nDjTnO5lzryJtAqS880wefmipwitVV7lLpAfpwwpeNPIERStgqF3EC76Bc/efOXpwWUtq6mFGcB5WQbIxYzk4Gc+nA0McbD41G7a0y6Y+EnGS+HU5N/DVivnsZTwamCIMb+twwo7lDMjApciCpGeGYQ7M6MJURiLg4LzraJqhZ64JSrWh5EYbxQJ5ITyKgbkGk22SixIddInrTnWtLJWBvgStebv3++SUBvpJM1tOnzR56GChxFyvF11H76r1dffvjqxYvb5z+/+fIHhRn0OumrAPg5RFN+mLo9pDxdb1Yo6GQsjXlhfTdXXzRRXKSOQhpw5tiP5R1KjxCI2mI2+Rv15lDy5MTSXkBfoaD3MOtecx3YziEyiwTQc4Gh2obZ3TMlqBye2fU/v/3d7a8/v/3k+ZsP/3Tz0R9vPv7k9cvfv3n2BxV1t0+f3/70mc7yq4VA6iJutKoJNOTziGuIaDRE7Ta++uAXRhByA1oOI5YZ6asPPs1NJhSVzUk2M7R0qSJFFdageWEhZ+RUSXUoyKUiqgFoixLWSBG6iIUL6maLvF8VU4m4k4xv9CbVmBEjhG7qm7CTIjgTvq+cFlEl9egdq27XNVdpI7YwxSxn+gM/nEKTkTDkA8QKpCD19QIZVq530B+2mv0MopQYRh6JhBuqD5S0lKKk9FgySlwtRwEuNH5u7IsnozXYOiudOZfGTP34ctb0RPrb/FRRQKB2HEWqIQr0U/44XkAznDww3c2Xr5a5/vdQ9ar5jjgb9zpFPoxzXCiERbYgSL5MTZql0lSPiXce1c4eCtNkQ+AVXjrBHPq1XnVal9BXRqbIO90oCqONySfjJLABdNKQUczM4ilGnqnqur3srBJCkvDSD3Xz0zsMO8QnvNxTtplJD8LLB8vC2ZcV7bGh72ljlrt6pGbJJsi/5wtH+X+f85WzNaMIGivxAFj6/T25uI7gvltSG7UVKShugLC6fiDfpxHjT0Qm0F+KcppSTlWH+kRcuNSHHJ86jUeKAAAzP7tZUDMJXa7+qhL6fTqnJa+FVhLfxUnAlouFauF65TxFlJZDcg0NZm4GSIrT70HZUO8oZWZTDgu5NJUXb5eMns94WcXemWLikalVdaekXt3yrO3qDiGNqm3D/9it1y3rkXj51JND5qDkquQiZnJg2V87nf8By4HADSYdAAA=")
